$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (Changed) date column C for rows 2-8 from 45224 (2023-10-25) to 45233 (2023-11-03)
foreach ($row in 2..8) {
    $ws.Cells.Item($row, 3).Value = 45233
}
